# Update benchmark: 2025-11-26 06:40:13 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ŞANS OYUNLARI
$ws.Range("D2").Value = '20.623,81 TL - 20.623,81 TL'
$ws.Range("H2").Value = '15 TL - 15 TL'

# Row 3 - HESAPTAN EFT - Şube (clear C3/D3)
$ws.Range("C3").Value = ''
$ws.Range("D3").Value = ''

# Row 4 - HESAPTAN EFT - ATM (clear C4/D4)
$ws.Range("C4").Value = ''
$ws.Range("D4").Value = ''

# Row 5 - HESAPTAN EFT - Mobil (clear C5/D5)
$ws.Range("C5").Value = ''
$ws.Range("D5").Value = ''

# Row 6 - DÜZENLİ EFT (clear C6/D6)
$ws.Range("C6").Value = ''
$ws.Range("D6").Value = ''

# Row 7 - KREDİ KARTINDAN FATURA ÖDEME
$ws.Range("H7").Value = '%3,09'

# Row 8 - HESAPTAN HAVALE - Şube (clear C8/D8)
$ws.Range("C8").Value = ''
$ws.Range("D8").Value = ''

# Row 9 - HESAPTAN HAVALE - ATM (clear C9/D9)
$ws.Range("C9").Value = ''
$ws.Range("D9").Value = ''

# Row 10 - HESAPTAN HAVALE - Mobil (clear C10/D10)
$ws.Range("C10").Value = ''
$ws.Range("D10").Value = ''

# Row 11 - DÜZENLİ HAVALE (clear C11/D11)
$ws.Range("C11").Value = ''
$ws.Range("D11").Value = ''

# Row 12 - GİDEN SWIFT
$ws.Range("C12").Value = 'WU: 0.29 USD–9,51 USD'
$ws.Range("D12").Value = 'WU: 0,75 USD–; Diğer: 909,5 TL–909,5 TL'

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = 'Hesaba: Asgari 15.714,29 TL | Azami 0,94 TL'
$ws.Range("D13").Value = 'Hesaba: Asgari 909,5 TL | Azami 909,5 TL'
$ws.Range("E13").Value = 'Hesaba: Asgari 1 TL | Azami 851,5 TL'
$ws.Range("K13").Value = 'Hesaba: Asgari 1 TL | Azami 865,75 TL'

# Row 14 - GİDEN SWIFT - Mobil (clear C14/D14)
$ws.Range("C14").Value = ''
$ws.Range("D14").Value = ''

# Row 15 - ÇEK TAHSİLİ BAŞKA BANKA
$ws.Range("D15").Value = '%0,8 Asgari Tutar:  Azami Tutar: '

# Row 17 - AYNI ŞUBE ÇEK TAHSİLATI
$ws.Range("D17").Value = '%0,8 Asgari Tutar:  Azami Tutar: '

# Row 20 - ÇEK İADE
$ws.Range("D20").Value = '20.623,81 TL'

# Row 21 - BLOKE ÇEK DÜZENLEME
$ws.Range("D21").Value = '%0,5 Asgari Tutar: 20.623,81 TL Azami Tutar: '

# Row 22 - YP ÇEK TAKASA GÖNDERME
$ws.Range("D22").Value = '%1 Asgari Tutar: 20.623,81 TL Azami Tutar: '

# Row 23 - ÇEK KARNESİ SAYFA ÜCRETİ
$ws.Range("D23").Value = '20.623,81 TL'

# Row 24 - SENET TAHSİLE ALMA
$ws.Range("D24").Value = '20.623,81 TL'

# Row 25 - MUAMELESİZ SENET İADESİ
$ws.Range("D25").Value = '20.623,81 TL'
